$d = $word.ActiveDocument
$p = $d.Paragraphs(1)
$xml = @"
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="Ttulo1"/>
      </w:pPr>
      <w:r>
        <w:t>RELATÓRIO DE DESENVOLVIMENTO</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve">: </w:t>
      </w:r>
      <w:r>
        <w:t>TRABALHO FINAL E</w:t>
      </w:r>
      <w:r>
        <w:t>STRUTURAS DE DADOS 2</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="SemEspaamento"/>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="SemEspaamento"/>
      </w:pPr>
      <w:r>
        <w:t xml:space="preserve">O trabalho proposto, foi desenvolvido utilizando como estrutura de dados principal a Árvore Binária de Busca, por ter sido a mais observada e testada em laboratório, além de ser mais didática de implementar. </w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve">Como forma de auxiliar o desenvolvimento do projeto, foi utilizado também a inteligência artificial </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>ChatGPT</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve"> para tirar dúvidas a respeito da implementação das funções, além de corrigir erros lógicos e sintáticos do código disponibilizado. </w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="SemEspaamento"/>
      </w:pPr>
      <w:r>
        <w:t>O desenvolvimento teve a seguinte sequência, uma função que interpretasse um arquivo .</w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>txt</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve"> que seria usado como entrada para o processamento do programa, ele deveria ser lido através do terminal</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> e suas informações armazenadas na estrutura. A estrutura tinha como propriedades: Valor, Valor convertido em binário, ponteiro para seu filho a esquerda e ponteiro para o seu filho a direita. Na atribuição das características do nó, já seria chamado a função </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>ConverterBinarioEmString</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve">, que recebia o valor atribuído ao nó, convertia ele para binário, adicionando zeros a esquerda se a quantidade de bits fosse menor do que o tamanho de N, e adicionava essa conversão na característica </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>BinarioString</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve"> do nó. Depois que a arvore inteira estivesse armazenada, com todos os seus valores convertidos para binário, seria chamado a função </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>ConcatenaStringArvore</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve">, que percorre a arvore inteira em ordem e concatenaria todos os valores binários de todos os nós retornando uma </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>string</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve"> armazenada no vetor </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>ArvoreConcatenada</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve">. A partir desse vetor, é chamado a função </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>ContarCombinacoesRepetidas</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve"> que cria as possíveis combinações de 0 e 1 de acordo com o tamanho da variável k, e </w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve">realiza a contagem que cada repetição apareceu na </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>string</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve"> fornecida. Além disso ela que realiza a impressão na tela das combinações e a quantidade de repetição de cada combinação.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="SemEspaamento"/>
      </w:pPr>
    </w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
"@
$p.Range.InsertXML($xml)
